$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column D ("Version Retired" -> shifts to E,
#    "Notes" -> shifts to F). This makes room for the new
#    "Version Deprecated" column.
$ws.Columns("D:D").Insert()

# 2. Set the new header for column D.
$ws.Range("D1").Value = "Version Deprecated"

# 3. Update row 73: "US Core Sex Extension" info is being replaced by the
#    "US Core Individual Sex Extension" entry, with cleared version/notes.
$ws.Range("B73").Value = "US Core Individual Sex Extension"
$ws.Range("C73").Value = "8.0.1"
$ws.Range("D73").Value = ""
$ws.Range("E73").Value = ""
$ws.Range("F73").Value = ""

# 4. Insert a new row before the final footnote row (row 77) to hold the
#    old "US Core Sex Extension" entry, now marked as superseded.
$ws.Rows("77:77").Insert()

$ws.Range("A77").Value = "     -"
$ws.Range("B77").Value = "US Core Sex Extension"
$ws.Range("C77").Value = "6.1.0"
$ws.Range("D77").Value = "8.0.1"
$ws.Range("E77").Value = ""
$ws.Range("F77").Value = "Superseded by US Core Individual Sex Extension in version 8.0.1"
